# Generate Report for Handback
# Adds a new handback record (f155536a-fc07-4934-84be-ea4978850230) as row 4
# on the "Overview", "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$guid = "f155536a-fc07-4934-84be-ea4978850230"
$hash = "bebd4163367f2c679c53b0fcc3158ace7e58a463"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f155536afc0749348 4beea4978850230/e2e/$guid.md",
    "",
    "",
    "$guid.md") | Out-Null
$wsOverview.Range("B4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C4").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f155536afc074934/e2e/$guid.md",
    "",
    "",
    "$guid.md") | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("B4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f155536afc074934/e2e/$guid.md",
    "",
    "",
    ".md") | Out-Null
$wsZh.Range("C4").Value = "Handed back: in sync with en-US"
$wsZh.Hyperlinks.Add(
    $wsZh.Range("D4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/$guid.$hash.zh-cn.xlf",
    "",
    "",
    "$guid.$hash.zh-cn.xlf") | Out-Null
$wsZh.Range("E4").Value = "2016-03-14 09:08:12"
$wsZh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f155536afc074934/e2e/$guid.md",
    "",
    "",
    "$guid.md") | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("G4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/$guid.$hash.zh-cn.xlf",
    "",
    "",
    "$guid.$hash.zh-cn.xlf") | Out-Null
$wsZh.Range("H4").Value = "2016-03-14 09:08:48"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = "Include"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f155536afc074934/e2e/$guid.md",
    "",
    "",
    "$guid.md") | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("B4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f155536afc074934/e2e/$guid.md",
    "",
    "",
    ".md") | Out-Null
$wsDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDe.Hyperlinks.Add(
    $wsDe.Range("D4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/$guid.$hash.de-de.xlf",
    "",
    "",
    "$guid.$hash.de-de.xlf") | Out-Null
$wsDe.Range("E4").Value = "2016-03-14 09:08:20"
$wsDe.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f155536afc074934/e2e/$guid.md",
    "",
    "",
    "$guid.md") | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("G4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/$guid.$hash.de-de.xlf",
    "",
    "",
    "$guid.$hash.de-de.xlf") | Out-Null
$wsDe.Range("H4").Value = "2016-03-14 09:09:01"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = "Include"

Write-Output "Handback report row added for $guid"
